$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted data (descending by value), matching the committed final state
$data = @(
    @("English", 21.06704746362185),
    @("Chinese", 20.25181588192941),
    @("Spanish", 6.231473634546495),
    @("Arabic", 4.146582678594735),
    @("German", 3.983758748423945),
    @("Japanese", 3.700624474687906),
    @("Malay-Indonesian", 3.226635296766861),
    @("Russian", 2.983797052721851),
    @("Portuguese", 2.753730537993464),
    @("French", 2.465060682794429),
    @("Turkish", 2.047644138663209),
    @("Italian", 1.831561013969273),
    @("Korean", 1.675942280999905),
    @("Dutch", 1.185570464338249),
    @("Polish", 0.9773998165646549),
    @("Persian", 0.969683360441615),
    @("Thai", 0.9121283300243744),
    @("Urdu", 0.907305759882559),
    @("Bengali", 0.8455053316711153),
    @("Vietnamese", 0.8315437217252539)
)

# Write the reordered data into rows 2..21
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (previously rows 22 and 23 - Uzbek & old Vietnamese)
$ws.Range("A22:B23").Delete()
